$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.934
$ws.Range("D6").Value = -7.88
$ws.Range("D7").Value = -7.411999999999999
$ws.Range("C8").Value = -12.672
$ws.Range("D8").Value = -7.858
$ws.Range("A12").Value = -21.864
$ws.Range("C12").Value = -13.002
$ws.Range("C14").Value = -11.675
$ws.Range("D19").Value = -7.640000000000001
$ws.Range("D21").Value = -7.742
$ws.Range("C22").Value = -12.929
$ws.Range("D24").Value = -7.821000000000001
